$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before the trailing summary block (old row 58 becomes row 60).
# Excel copies formatting down from the row above (row 57), which is exactly what
# we want for the two new data rows (58 and 59).
$ws.Rows("58:60").Insert()

# Row 58 -> new test case DEC_0871 (mirrors the pattern of rows 2-57)
$ws.Range("A58").Value = "DEC_0871"
$ws.Range("B58").Value = "18092588-0"
$ws.Range("C58").Value = "'sebA`$1357"
$ws.Range("D58:J58").Value = "SIN_DATO"

# Row 59 -> new test case DEC_0872
$ws.Range("A59").Value = "DEC_0872"
$ws.Range("B59").Value = "18092588-0"
$ws.Range("C59").Value = "'sebA`$1357"
$ws.Range("D59:J59").Value = "SIN_DATO"

# Row 60 is the blank spacer row that used to be row 58 (only B/C formatted, no
# other content) -- clear the stray cells that the row-insert operation created
# in columns A and D:J so it matches the original spacer row exactly.
$ws.Range("A60").Clear()
$ws.Range("D60:J60").Clear()

$ws.Range("A59").Select()
